# Re-process the sheet with the newly curated dimensions:
#  - "lugar-de-residencia" and "lugar-de-nacimiento" move from
#    iaest-dimension:* to iaest-measure:* (row 2, columns A & D)
#  - row 3/4 "dim"/"medida"/"xsd:int" roles are re-assigned for
#    columns A, C, D, F, G accordingly
#  - row 5 (the mapping-*.xlsx helper row) is no longer needed and removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 — dimension/measure qualifiers
$ws.Range("A2").Value = "iaest-measure:lugar-de-residencia"
$ws.Range("D2").Value = "iaest-measure:lugar-de-nacimiento"

# Row 3 — dim / medida roles
$ws.Range("A3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "dim"

# Row 4 — data types / URIs
$ws.Range("A4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"

# Row 5 no longer exists — clear it entirely
$ws.Range("A5:G5").Clear()
